$d = $word.ActiveDocument

# Locate the paragraph that currently reads:
#   "IF obstacle is to the right move to the left"
$searchRange = $d.Content
$found = $searchRange.Find.Execute("IF obstacle is to the right move to the left", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target paragraph text."
}

$targetPara = $searchRange.Paragraphs(1)
$paraRange = $targetPara.Range
$paraEnd = $paraRange.End

# Append " and vice versa." to the end of that paragraph's text (before the
# paragraph mark).
$insertPoint = $d.Range($paraEnd - 1, $paraEnd - 1)
$insertPoint.InsertAfter(" and vice versa.")

# Insert a brand new list paragraph right after it; InsertParagraphAfter
# naturally inherits the ListParagraph style / numPr (ilvl 3, numId 3) of
# the paragraph it is split from.
$updatedParaRange = $targetPara.Range
$endPoint = $d.Range($updatedParaRange.End - 1, $updatedParaRange.End - 1)
$endPoint.InsertParagraphAfter()

$newPara = $targetPara.Next()
$newRange = $newPara.Range
$newRange.InsertAfter("Seemingly each engine may only defuzzify once and that is it so I am unsure how to have the engines work around this. For now I will use multiple engines.")
